$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are stored as text (e.g. "60.220.25"); prefix numeric-looking
# ones with an apostrophe so Excel keeps them as text instead of coercing to a
# number, then reset the style so no stray number-format style is introduced.
$ws.Range("D2").Value = "'60.220.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "'2.601.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'582.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "'143.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "'6.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "'3.059.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "'24.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "'60.219.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "'0.0000140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'2.605.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'11.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'345.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'6.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "'63.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "'8.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").Value = "'1.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.43%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'6.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'167.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  +9.92%  "
$ws.Range("D35").Value = "'4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "'0.984"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").Value = "'38.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "'312.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").Value = "'3.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'0.844"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").Value = "'135.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").Value = "'0.0995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'19.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "'0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'0.0243"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").Value = "'19.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("E51").Value = "  +0.51%  "
